# Actualización automática 2025-09-08 14:45:08
$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO": update PORCELANATO total for row 8 (client FRANK FERRETERIA FRANKFERRE CIA.)
# and bump the "x de 8" progress counter in M10 from "1 de 8" to "2 de 8".
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M8").Value = 5372.02
$wsGrupo.Range("M10").Value = "2 de 8"

# Sheet "VENTA MENSUAL": update "septiembre" sales for the same client (F8)
# and recompute the PRESUPUESTO-period total in F10 accordingly.
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F8").Value = 5372.02
$wsMensual.Range("F10").Value = 5428.400000000001
